# Decrement the "Moose_new" (column S) value by 1 for the moose-level
# records affected by this edit (rows 142-345), skipping the block of
# rows (198-201) that were already at the correct level and therefore
# untouched by the original commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$skipRows = @(198, 199, 200, 201)

for ($row = 142; $row -le 345; $row++) {
    if ($skipRows -contains $row) {
        continue
    }

    $cell = $ws.Cells.Item($row, 19)   # column S = 19
    $current = $cell.Value2
    $cell.Value2 = $current - 1
}
